# Generate Report for handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   for the two localized entries (rows 2 & 3) in both the zh-cn and de-de sheets.
# - The "Latest Target File" (E) and "Latest Handback File" (F) columns get
#   populated with hyperlinks to the source markdown file and the handed-back
#   xlf file, respectively.
# - The "Latest Handback DateTime" (G) column is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/438f75c76e5b9cf6c20f68363c2ab4cdcf6fa388/e2e/8a486b4e-032b-4849-b688-a223ba0c2641.md"
$mdName  = "8a486b4e-032b-4849-b688-a223ba0c2641.md"

$statusHandedBack = "Handed back: in sync with en-US"

function Update-LangSheet {
    param([string]$sheetName, [string]$xlfUrl, [string]$xlfName, [string]$handbackTime)

    $ws = $wb.Worksheets($sheetName)

    # Row 2 and row 3 both refer to the same handed-off source + target xlf.
    foreach ($row in 2, 3) {
        $ws.Range("B$row").Value = $statusHandedBack

        $eCell = $ws.Range("E$row")
        $ws.Hyperlinks.Add($eCell, $mdUrl, "", "", $mdName) | Out-Null

        $fCell = $ws.Range("F$row")
        $ws.Hyperlinks.Add($fCell, $xlfUrl, "", "", $xlfName) | Out-Null

        $ws.Range("G$row").Value = $handbackTime
    }
}

Update-LangSheet "zh-cn" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6608ae60a853fddd445a1f998697db534646a6ab/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/8a486b4e-032b-4849-b688-a223ba0c2641.c77709331ccd18e339cb2be3d0e883b76063674d.zh-cn.xlf" `
    "8a486b4e-032b-4849-b688-a223ba0c2641.c77709331ccd18e339cb2be3d0e883b76063674d.zh-cn.xlf" `
    "2016-01-25 09:11:26"

Update-LangSheet "de-de" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/48b89a48914a0876a22bd72e30406de663e3d337/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/8a486b4e-032b-4849-b688-a223ba0c2641.c77709331ccd18e339cb2be3d0e883b76063674d.de-de.xlf" `
    "8a486b4e-032b-4849-b688-a223ba0c2641.c77709331ccd18e339cb2be3d0e883b76063674d.de-de.xlf" `
    "2016-01-25 09:11:44"

# The Overview sheet mirrors the same "Status" shared string for each
# per-language column, so it reflects the handback too.
$overview = $wb.Worksheets("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

Write-Output "Report generated for handback."
